# Adding reseller test & data
#
# 1. Insert a new worksheet "searchResellerInvalidID" right after "TestCases"
#    (i.e. as the second tab, before "searchResellerID"), populated with a
#    reseller-search test case that looks up an id that doesn't exist and
#    expects no results.
# 2. Record the new test case on the "TestCases" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Register the new test case on the TestCases sheet ---------------
$testCases = $wb.Worksheets.Item("TestCases")
$testCases.Cells.Item(6, 1).Value = "searchResellerInvalidID"
$testCases.Cells.Item(6, 3).Value = "Y"

# --- 2. Insert the new sheet right after TestCases -----------------------
$searchResellerID = $wb.Worksheets.Item("searchResellerID")
$newSheet = $wb.Worksheets.Add($searchResellerID)
$newSheet.Name = "searchResellerInvalidID"

# Header row (same shape as the other searchReseller* sheets, plus a new
# "expected" column describing the anticipated outcome).
$headers = @("userName", "password", "resName", "id", "searchList", "expected")
for ($col = 1; $col -le $headers.Length; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$newSheet.Range("A1:F1").Interior.Color = 65535

# Data row: same credentials/resName as the other reseller searches, but an
# id that doesn't exist, so the query is expected to return nothing.
$newSheet.Cells.Item(2, 1).Value = "admin"
$newSheet.Cells.Item(2, 2).Value = "123qwe"
$newSheet.Cells.Item(2, 3).Value = "autoTestReseller1"
$newSheet.Cells.Item(2, 4).Value = -9444
$newSheet.Cells.Item(2, 5).Value = "Reseller ID"
$newSheet.Cells.Item(2, 6).Value = "Query returned no results."

# --- 3. Restore cursor positions on the touched sheets --------------------
$testCases.Range("D10").Select()
$newSheet.Range("G9").Select()
